$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Won" (column D) and "Percent" (column E) values per tournament round,
# recomputed against the bracket/prediction data. Percent is stored as a
# literal text string (e.g. "75%"), so force text entry and then restore
# the default "Normal" style so no new cell formatting is introduced.
$ws.Range("E2:E9").NumberFormat = "@"

$ws.Range("D2").Value = 3
$ws.Range("E2").Value = "75%"

$ws.Range("D3").Value = 20
$ws.Range("E3").Value = "62%"

$ws.Range("D4").Value = 11
$ws.Range("E4").Value = "69%"

$ws.Range("D5").Value = 3
$ws.Range("E5").Value = "38%"

$ws.Range("D6").Value = 2
$ws.Range("E6").Value = "50%"

$ws.Range("D7").Value = 1
$ws.Range("E7").Value = "50%"

$ws.Range("D8").Value = 1
$ws.Range("E8").Value = "100%"

$ws.Range("D9").Value = 41
$ws.Range("E9").Value = "61%"

$ws.Range("E2:E9").Style = "Normal"
